$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sản phẩm")

# --- Column B (Tên MH / product name) -------------------------------------
# Set values first (in row order) so the shared-string table gets the new
# names appended in this order, then strip the old wrap-text style so the
# cells fall back to the default style (matches target: no s="2").
$ws.Range("B2").Value = "AMD Ryzen 7 5800X3D"
$ws.Range("B3").Value = "GIGABYTE GeForce RTX 4060 Ti EAGLE 8G"
$ws.Range("B4").Value = "Corsair Vengeance RS RGB 2x16GB 3600 "
$ws.Range("B5").Value = "Asus X570 Rog Crosshair VIII Dark Hero"
$ws.Range("B6").Value = "ASUS TUF Gaming 750B - 80 Plus Bronze (750W)"
$ws.Range("B7").Value = "SSD WD Black SN850x 1TB M.2 PCIe NVMe Gen 4.0"
$ws.Range("B8").Value = "Màn hình cong LG 27GS60QC-B UltraGear 27`" 2K 180Hz chuyên game"
$ws.Range("B9").Value = "NZXT H5 Flow Black"
$ws.Range("B10").Value = "Bộ 3 quạt Corsair RS120 ARGB BLACK (CO-9050181-WW)"
$ws.Range("B11").Value = "Cooler Master HYPER 620S"

$ws.Range("B2:B10").ClearFormats()

# --- Column E (Loại sản phẩm / category) -----------------------------------
$ws.Range("E2").Value = "Processor"
$ws.Range("E3").Value = "Graphics Card"
$ws.Range("E4").Value = "RAM"
$ws.Range("E5").Value = "Motherboard"
$ws.Range("E6").Value = "PSU"
$ws.Range("E7").Value = "SSD"
$ws.Range("E8").Value = "Monitor"
$ws.Range("E9").Value = "Case"
$ws.Range("E10").Value = "Cooling"
$ws.Range("E11").Value = "Cooling"

# --- Column F (Số lượng cần mua / qty) -> all become 1 ----------------------
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = 1

# --- New row 11: Barcode + unit, matching the formatting of the rows above -
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = "123456789010"

$ws.Range("D11").Value = "Chiếc"

$ws.Range("F10").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("F11").Value = 1

$excel.CutCopyMode = $false

# --- Extend Table1 to cover the new row ------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F11"))

# --- Restore the cursor position reflected in the saved file ---------------
$ws.Range("J7").Select() | Out-Null
